# Append two new log rows (138 and 139) to the feed logs sheet,
# matching the existing run_id / rss_url_id / date / response / item_count layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 138
$ws.Cells.Item(138, 1).Value = 137
$ws.Cells.Item(138, 2).Value = 1
$ws.Cells.Item(138, 3).Value = "2024-06-17 20:11:19"
$ws.Cells.Item(138, 4).Value = 200
$ws.Cells.Item(138, 5).Value = 14

# Row 139
$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = 2
$ws.Cells.Item(139, 3).Value = "2024-06-17 20:11:20"
$ws.Cells.Item(139, 4).Value = 200
$ws.Cells.Item(139, 5).Value = 3
